$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "tcs"
$ws.Range("B2").Value = "tata"
$ws.Range("C2").Value = "services"

# Delete rows 3-6 (old WPS2..Pdf5 rows)
$ws.Range("A3:C6").EntireRow.Delete()

# Update selection to match target state
$ws.Range("C5").Select()
